$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to remain text so numeric-looking strings
# (e.g. trailing zeros, multi-dot formatted numbers) are preserved exactly.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '71.655.32'
$ws.Range('E2').Value = '  -0.62%  '
$ws.Range('D3').Value = '3.985.14'
$ws.Range('E3').Value = '  -0.93%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '541.18'
$ws.Range('E5').Value = '  +4.87%  '
$ws.Range('D6').Value = '150.16'
$ws.Range('E6').Value = '  +2.04%  '
$ws.Range('D7').Value = '0.702'
$ws.Range('E7').Value = '  +2.13%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '0.746'
$ws.Range('E9').Value = '  -1.40%  '
$ws.Range('E10').Value = '  -2.99%  '
$ws.Range('D11').Value = '53.66'
$ws.Range('E11').Value = '  +14.63%  '
$ws.Range('D12').Value = '0.0000323'
$ws.Range('E12').Value = '  -1.13%  '
$ws.Range('D13').Value = '10.65'
$ws.Range('E13').Value = '  -1.91%  '
$ws.Range('D14').Value = '4.629.00'
$ws.Range('E14').Value = '  -1.01%  '
$ws.Range('D15').Value = '3.992.10'
$ws.Range('E15').Value = '  -0.97%  '
$ws.Range('D16').Value = '14.15'
$ws.Range('E16').Value = '  -0.11%  '
$ws.Range('D17').Value = '20.51'
$ws.Range('E17').Value = '  -3.33%  '
$ws.Range('E18').Value = '  -0.33%  '
$ws.Range('D19').Value = '1.18'
$ws.Range('E19').Value = '  -1.71%  '
$ws.Range('D20').Value = '71.707.73'
$ws.Range('E20').Value = '  -0.54%  '
$ws.Range('D21').Value = '430.29'
$ws.Range('E21').Value = '  -1.45%  '
$ws.Range('B22').Value = 'ImmutableX'
$ws.Range('C22').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D22').Value = '3.57'
$ws.Range('E22').Value = '  +0.45%  '
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').Value = '96.64'
$ws.Range('E23').Value = '  -4.83%  '
$ws.Range('D24').Value = '4.28'
$ws.Range('E24').Value = '  +7.19%  '
$ws.Range('E25').Value = '  -2.73%  '
$ws.Range('D26').Value = '11.46'
$ws.Range('E26').Value = '  -2.22%  '
$ws.Range('D27').Value = '10.64'
$ws.Range('E27').Value = '  -4.86%  '
$ws.Range('D28').Value = '5.86'
$ws.Range('E28').Value = '  +1.13%  '
$ws.Range('D29').Value = '36.71'
$ws.Range('E29').Value = '  -2.07%  '
$ws.Range('E30').Value = '  +18.70%  '
$ws.Range('D31').Value = '7.48'
$ws.Range('E31').Value = '  +8.39%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = '0.132'
$ws.Range('E32').Value = '  +2.61%  '
$ws.Range('B33').Value = 'Cosmos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D33').Value = '13.42'
$ws.Range('E33').Value = '  -0.74%  '
$ws.Range('D34').Value = '49.04'
$ws.Range('E34').Value = '  +17.86%  '
$ws.Range('D35').Value = '678.51'
$ws.Range('E35').Value = '  -0.80%  '
$ws.Range('D36').Value = '65.93'
$ws.Range('E36').Value = '  -3.04%  '
$ws.Range('D37').Value = '0.441'
$ws.Range('E37').Value = '  +0.51%  '
$ws.Range('D38').Value = '0.0₃0830'
$ws.Range('E38').Value = '  -5.08%  '
$ws.Range('D39').Value = '0.151'
$ws.Range('E39').Value = '  -0.38%  '
$ws.Range('E40').Value = '  -4.92%  '
$ws.Range('D41').Value = '3.37'
$ws.Range('E41').Value = '  +3.80%  '
$ws.Range('D42').Value = '0.999'
$ws.Range('E42').Value = '  +0.19%  '
$ws.Range('E43').Value = '  +0.29%  '
$ws.Range('E44').Value = '  -0.63%  '
$ws.Range('E45').Value = '  +0.39%  '
$ws.Range('D46').Value = '0.149'
$ws.Range('E46').Value = '  -4.38%  '
$ws.Range('E47').Value = '  +9.36%  '
$ws.Range('D48').Value = '3.38'
$ws.Range('E48').Value = '  -3.55%  '
$ws.Range('D49').Value = '0.000280'
$ws.Range('E49').Value = '  +3.86%  '
$ws.Range('D50').Value = '3.00'
$ws.Range('E50').Value = '  -2.70%  '
$ws.Range('D51').Value = '144.57'
$ws.Range('E51').Value = '  +1.48%  '

# Restore the original (default) number format/style on column D so the
# cell styling matches the source workbook.
$ws.Range('D2:D51').NumberFormat = 'General'
$ws.Range('D2:D51').Style = 'Normal'
